$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 118: Crafty Concoctions | Commanding Craftsman's Syrup
$ws.Cells.Item(118, 8).Value = 1692.125
$ws.Cells.Item(118, 9).Value = 950.44446
$ws.Cells.Item(118, 10).Value = 2645.7144
$ws.Cells.Item(118, 11).Value = 2851.33338
$ws.Cells.Item(118, 12).Value = 7937.1432
$ws.Cells.Item(118, 13).Value = -1194.33338
$ws.Cells.Item(118, 14).Value = -11251.1432

# Row 125: Body over Mind | Grade 5 Dexterity Alkahest
$ws.Cells.Item(125, 8).Value = 20081.416
$ws.Cells.Item(125, 9).Value = 3576
$ws.Cells.Item(125, 10).Value = 36586.832
$ws.Cells.Item(125, 11).Value = 32184
$ws.Cells.Item(125, 12).Value = 329281.488
$ws.Cells.Item(125, 13).Value = -29724
$ws.Cells.Item(125, 14).Value = -334201.488

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Cells.Item(45, 8).Value = 995
$ws.Cells.Item(45, 9).Value = 995
$ws.Cells.Item(45, 11).Value = 995
$ws.Cells.Item(45, 13).Value = -618

# Row 63: Rivets Run through It | Mythrite Rivets
$ws.Cells.Item(63, 8).Value = 9999.857
$ws.Cells.Item(63, 10).Value = 9999.857
$ws.Cells.Item(63, 12).Value = 9999.857
$ws.Cells.Item(63, 14).Value = -11371.857

# Row 66: A Riveting Revival (L) | Mythrite Rivets
$ws.Cells.Item(66, 8).Value = 9999.857
$ws.Cells.Item(66, 10).Value = 9999.857
$ws.Cells.Item(66, 12).Value = 49999.285
$ws.Cells.Item(66, 14).Value = -56863.285

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Cells.Item(74, 8).Value = 1937.5
$ws.Cells.Item(74, 9).Value = 1850.6154
$ws.Cells.Item(74, 10).Value = 2314
$ws.Cells.Item(74, 11).Value = 1850.6154
$ws.Cells.Item(74, 12).Value = 2314
$ws.Cells.Item(74, 13).Value = -976.6153999999999
$ws.Cells.Item(74, 14).Value = -4062

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Cells.Item(77, 8).Value = 1937.5
$ws.Cells.Item(77, 9).Value = 1850.6154
$ws.Cells.Item(77, 10).Value = 2314
$ws.Cells.Item(77, 11).Value = 9253.076999999999
$ws.Cells.Item(77, 12).Value = 11570
$ws.Cells.Item(77, 13).Value = -4885.076999999999
$ws.Cells.Item(77, 14).Value = -20306

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Cells.Item(122, 8).Value = 12349618
$ws.Cells.Item(122, 9).Value = 27781316
$ws.Cells.Item(122, 11).Value = 83343948
$ws.Cells.Item(122, 13).Value = -83341498

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Cells.Item(99, 8).Value = 2359.7778
$ws.Cells.Item(99, 9).Value = 1262
$ws.Cells.Item(99, 11).Value = 1262
$ws.Cells.Item(99, 13).Value = 236

# Row 107: The Gold Experience | Deepgold Nugget
$ws.Cells.Item(107, 8).Value = 1311.5625
$ws.Cells.Item(107, 9).Value = 1332.3334
$ws.Cells.Item(107, 10).Value = 1000
$ws.Cells.Item(107, 11).Value = 1332.3334
$ws.Cells.Item(107, 12).Value = 1000
$ws.Cells.Item(107, 13).Value = 587.6666
$ws.Cells.Item(107, 14).Value = -4840

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Cells.Item(31, 8).Value = 6668.091
$ws.Cells.Item(31, 9).Value = 2703.8333
$ws.Cells.Item(31, 10).Value = 11425.2
$ws.Cells.Item(31, 11).Value = 2703.8333
$ws.Cells.Item(31, 12).Value = 11425.2
$ws.Cells.Item(31, 13).Value = -2408.8333
$ws.Cells.Item(31, 14).Value = -12015.2

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Cells.Item(34, 8).Value = 6668.091
$ws.Cells.Item(34, 9).Value = 2703.8333
$ws.Cells.Item(34, 10).Value = 11425.2
$ws.Cells.Item(34, 11).Value = 2703.8333
$ws.Cells.Item(34, 12).Value = 11425.2
$ws.Cells.Item(34, 13).Value = -2501.8333
$ws.Cells.Item(34, 14).Value = -11829.2

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Cells.Item(58, 8).Value = 2592.7273
$ws.Cells.Item(58, 9).Value = 2389
$ws.Cells.Item(58, 10).Value = 2949.25
$ws.Cells.Item(58, 11).Value = 2389
$ws.Cells.Item(58, 12).Value = 2949.25
$ws.Cells.Item(58, 13).Value = -2186
$ws.Cells.Item(58, 14).Value = -3355.25

# Row 62: Splinter in the Sewers | Cedar Lumber
$ws.Cells.Item(62, 8).Value = 5729.8335
$ws.Cells.Item(62, 10).Value = 6221.25
$ws.Cells.Item(62, 12).Value = 6221.25
$ws.Cells.Item(62, 14).Value = -7469.25

# Row 65: The Lumber of Their Discontent (L) | Cedar Lumber
$ws.Cells.Item(65, 8).Value = 5729.8335
$ws.Cells.Item(65, 10).Value = 6221.25
$ws.Cells.Item(65, 12).Value = 31106.25
$ws.Cells.Item(65, 14).Value = -37346.25

# Row 99: O Pine | Pine Lumber
$ws.Cells.Item(99, 8).Value = 2793.6667
$ws.Cells.Item(99, 9).Value = 2352.4
$ws.Cells.Item(99, 11).Value = 2352.4
$ws.Cells.Item(99, 13).Value = -854.4000000000001

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Cells.Item(126, 8).Value = 2793.6667
$ws.Cells.Item(126, 9).Value = 2352.4
$ws.Cells.Item(126, 11).Value = 7057.200000000001
$ws.Cells.Item(126, 13).Value = -4587.200000000001

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Cells.Item(132, 8).Value = 2002.6086
$ws.Cells.Item(132, 9).Value = 1956.2778
$ws.Cells.Item(132, 10).Value = 2169.4
$ws.Cells.Item(132, 11).Value = 5868.8334
$ws.Cells.Item(132, 12).Value = 6508.200000000001
$ws.Cells.Item(132, 13).Value = -3338.8334
$ws.Cells.Item(132, 14).Value = -11568.2

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Cells.Item(136, 8).Value = 2592.7273
$ws.Cells.Item(136, 9).Value = 2389
$ws.Cells.Item(136, 10).Value = 2949.25
$ws.Cells.Item(136, 11).Value = 7167
$ws.Cells.Item(136, 12).Value = 8847.75
$ws.Cells.Item(136, 13).Value = -4617
$ws.Cells.Item(136, 14).Value = -13947.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 18: Fisher of Men | Salt Cod
$ws.Cells.Item(18, 8).Value = 2241.0908
$ws.Cells.Item(18, 9).Value = 348.91666
$ws.Cells.Item(18, 10).Value = 4511.7
$ws.Cells.Item(18, 11).Value = 1046.74998
$ws.Cells.Item(18, 12).Value = 13535.1
$ws.Cells.Item(18, 13).Value = -877.7499800000001
$ws.Cells.Item(18, 14).Value = -13873.1

# Row 33: Cooking with Gas | Chicken Stock
$ws.Cells.Item(33, 8).Value = 175.66667
$ws.Cells.Item(33, 10).Value = 50
$ws.Cells.Item(33, 12).Value = 300
$ws.Cells.Item(33, 14).Value = -866

# Row 34: Fever Pitch | Chamomile Tea
$ws.Cells.Item(34, 8).Value = 275.8
$ws.Cells.Item(34, 10).Value = 299.5
$ws.Cells.Item(34, 12).Value = 898.5
$ws.Cells.Item(34, 14).Value = -1066.5

# Row 49: Leek Soup for the Soul | Cawl Cennin
$ws.Cells.Item(49, 8).Value = 349.5
$ws.Cells.Item(49, 9).Value = 399.66666
$ws.Cells.Item(49, 11).Value = 1198.99998
$ws.Cells.Item(49, 13).Value = -1042.99998

# Row 107: Slippery Service | Frantoio Oil
$ws.Cells.Item(107, 8).Value = 3230.4482
$ws.Cells.Item(107, 10).Value = 2303.6365
$ws.Cells.Item(107, 12).Value = 6910.9095
$ws.Cells.Item(107, 14).Value = -10750.9095

# Row 109: Cure for What Ails | Purple Carrot Juice
$ws.Cells.Item(109, 8).Value = 923761.4
$ws.Cells.Item(109, 9).Value = 1255121.9
$ws.Cells.Item(109, 10).Value = 40133.332
$ws.Cells.Item(109, 11).Value = 3765365.7
$ws.Cells.Item(109, 12).Value = 120399.996
$ws.Cells.Item(109, 13).Value = -3764325.7
$ws.Cells.Item(109, 14).Value = -122479.996

# Row 121: A Cookie for Your Troubles | Coffee Biscuit
$ws.Cells.Item(121, 8).Value = 849.3333
$ws.Cells.Item(121, 9).Value = 756
$ws.Cells.Item(121, 10).Value = 1036
$ws.Cells.Item(121, 11).Value = 2268
$ws.Cells.Item(121, 12).Value = 3108
$ws.Cells.Item(121, 13).Value = -958
$ws.Cells.Item(121, 14).Value = -5728

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Cells.Item(131, 8).Value = 2779299.8
$ws.Cells.Item(131, 9).Value = 612.0526
$ws.Cells.Item(131, 11).Value = 1836.1578
$ws.Cells.Item(131, 13).Value = 3203.8422

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5: Hora at Me | Bone Hora
$ws.Cells.Item(5, 8).Value = 1000
$ws.Cells.Item(5, 9).Value = 1000
$ws.Cells.Item(5, 11).Value = 1000
$ws.Cells.Item(5, 13).Value = -888

# Row 11: A Ringing Success | Copper Ring
$ws.Cells.Item(11, 8).Value = 10125934
$ws.Cells.Item(11, 10).Value = 51749.5
$ws.Cells.Item(11, 12).Value = 51749.5
$ws.Cells.Item(11, 14).Value = -52027.5

# Row 12: Horn of Plenty | Bone Armillae
$ws.Cells.Item(12, 8).Value = 257666.88
$ws.Cells.Item(12, 9).Value = 181431.7
$ws.Cells.Item(12, 10).Value = 473666.5
$ws.Cells.Item(12, 11).Value = 181431.7
$ws.Cells.Item(12, 12).Value = 473666.5
$ws.Cells.Item(12, 13).Value = -181291.7
$ws.Cells.Item(12, 14).Value = -473946.5

# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Cells.Item(70, 8).Value = 10110549
$ws.Cells.Item(70, 9).Value = 333333340
$ws.Cells.Item(70, 10).Value = 9837.125
$ws.Cells.Item(70, 11).Value = 333333340
$ws.Cells.Item(70, 12).Value = 9837.125
$ws.Cells.Item(70, 13).Value = -333333070
$ws.Cells.Item(70, 14).Value = -10377.125

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Cells.Item(73, 8).Value = 10110549
$ws.Cells.Item(73, 9).Value = 333333340
$ws.Cells.Item(73, 10).Value = 9837.125
$ws.Cells.Item(73, 11).Value = 333333340
$ws.Cells.Item(73, 12).Value = 9837.125
$ws.Cells.Item(73, 13).Value = -333332404
$ws.Cells.Item(73, 14).Value = -11709.125

# Row 86: Keeping Claw and Order | Griffin Talon Ring of Aiming
$ws.Cells.Item(86, 8).Value = 55000
$ws.Cells.Item(86, 10).Value = 55000
$ws.Cells.Item(86, 12).Value = 55000
$ws.Cells.Item(86, 14).Value = -57372

# Row 89: Ring of Reciprocity (L) | Griffin Talon Ring of Aiming
$ws.Cells.Item(89, 8).Value = 55000
$ws.Cells.Item(89, 10).Value = 55000
$ws.Cells.Item(89, 12).Value = 165000
$ws.Cells.Item(89, 14).Value = -176856

# Row 92: Play It by Ear | Triphane Earrings of Healing
$ws.Cells.Item(92, 8).Value = 21250
$ws.Cells.Item(92, 10).Value = 21250
$ws.Cells.Item(92, 12).Value = 21250
$ws.Cells.Item(92, 14).Value = -24994

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Cells.Item(102, 8).Value = 3392
$ws.Cells.Item(102, 9).Value = 3392
$ws.Cells.Item(102, 11).Value = 3392
$ws.Cells.Item(102, 13).Value = -1770

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Cells.Item(122, 8).Value = 9242.125
$ws.Cells.Item(122, 9).Value = 11852.63
$ws.Cells.Item(122, 10).Value = 3820.3076
$ws.Cells.Item(122, 11).Value = 35557.89
$ws.Cells.Item(122, 12).Value = 11460.9228
$ws.Cells.Item(122, 13).Value = -33107.89
$ws.Cells.Item(122, 14).Value = -16360.9228

# Row 132: On Board for Lar | Lar Ingot
$ws.Cells.Item(132, 8).Value = 3178.76
$ws.Cells.Item(132, 9).Value = 3146.1428
$ws.Cells.Item(132, 10).Value = 3350
$ws.Cells.Item(132, 11).Value = 9438.428400000001
$ws.Cells.Item(132, 12).Value = 10050
$ws.Cells.Item(132, 13).Value = -6908.428400000001
$ws.Cells.Item(132, 14).Value = -15110

# Row 133: Pendulums of Our Own | Lar Pendulums
$ws.Cells.Item(133, 8).Value = 100599.4
$ws.Cells.Item(133, 10).Value = 100599.4
$ws.Cells.Item(133, 12).Value = 100599.4
$ws.Cells.Item(133, 14).Value = -110719.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Cells.Item(7, 8).Value = 6700
$ws.Cells.Item(7, 9).Value = 6700
$ws.Cells.Item(7, 11).Value = 6700
$ws.Cells.Item(7, 13).Value = -6588

# Row 20: Choke Hold | Hard Leather Choker
$ws.Cells.Item(20, 8).Value = 308
$ws.Cells.Item(20, 9).Value = 1500
$ws.Cells.Item(20, 10).Value = 10
$ws.Cells.Item(20, 11).Value = 1500
$ws.Cells.Item(20, 12).Value = 10
$ws.Cells.Item(20, 13).Value = -1274
$ws.Cells.Item(20, 14).Value = -462

# Row 55: It's Not a Job, It's a Calling | Peiste Leather
$ws.Cells.Item(55, 8).Value = 1545.75
$ws.Cells.Item(55, 9).Value = 227.5
$ws.Cells.Item(55, 11).Value = 227.5
$ws.Cells.Item(55, 13).Value = -54.5

# Row 61: Spelling Me Softly | Raptor Leather
$ws.Cells.Item(61, 8).Value = 4740.35
$ws.Cells.Item(61, 9).Value = 2380.2
$ws.Cells.Item(61, 11).Value = 2380.2
$ws.Cells.Item(61, 13).Value = -2178.2

# Row 113: Peace in Rest | Atrociraptor Leather
$ws.Cells.Item(113, 8).Value = 4740.35
$ws.Cells.Item(113, 9).Value = 2380.2
$ws.Cells.Item(113, 11).Value = 2380.2
$ws.Cells.Item(113, 13).Value = -210.1999999999998

# Row 126: Battered Books | Saiga Leather
$ws.Cells.Item(126, 8).Value = 6700
$ws.Cells.Item(126, 9).Value = 6700
$ws.Cells.Item(126, 11).Value = 20100
$ws.Cells.Item(126, 13).Value = -17630

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 45: Private Concerns | Linen Trousers
$ws.Cells.Item(45, 8).Value = 10000
$ws.Cells.Item(45, 9).Value = 10000
$ws.Cells.Item(45, 11).Value = 10000
$ws.Cells.Item(45, 13).Value = -9509

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Cells.Item(136, 8).Value = 2597.111
$ws.Cells.Item(136, 10).Value = 3199.6667
$ws.Cells.Item(136, 12).Value = 9599.000100000001
$ws.Cells.Item(136, 14).Value = -14699.0001
